# atualização preços 25/09 e contrato 039 emporia
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Update the "futuros" price curve (source of most recalculated
#    values elsewhere in the workbook).
# ---------------------------------------------------------------
$futuros = $wb.Worksheets.Item("futuros")
$futuros.Range("D1").Value = 45925

$futuros.Range("B2").Value  = 368.6
$futuros.Range("B3").Value  = 349.1
$futuros.Range("B4").Value  = 335.85
$futuros.Range("B5").Value  = 322.7
$futuros.Range("B6").Value  = 310.25
$futuros.Range("B7").Value  = 301.5
$futuros.Range("B8").Value  = 298.8
$futuros.Range("B9").Value  = 292.4
$futuros.Range("B10").Value = 287.3
$futuros.Range("B11").Value = 281.6
$futuros.Range("B12").Value = 276.25

# ---------------------------------------------------------------
# 2) "Sheet2" - add pricing for contract row 36 (novo contrato 039
#    Emporia) and move the selection to P32.
# ---------------------------------------------------------------
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Range("O36").Value = 5.3368
$sheet2.Range("P36").Formula = "=M36*O36"
$sheet2.Range("Q36").Formula = "=P36*E36"

$sheet2.Range("P32").Select()

# ---------------------------------------------------------------
# 3) "hedge" becomes the active sheet/tab, with L16:L28 selected.
# ---------------------------------------------------------------
$hedge = $wb.Worksheets.Item("hedge")
$hedge.Activate()
$hedge.Range("A2").Select()
$hedge.Range("L16:L28").Select()
